$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("People")

# The previously-empty placeholder row 36 is no longer needed once row 35
# is filled in below, so its (empty, formatted-only) cells are cleared out.
$ws1.Rows("36").Clear()

# Fill in the new "Martin Luther University Halle-Wittenberg" / "DiP-KS"
# entry on row 35 (institution columns first, then the name/link columns,
# matching the order the strings were actually entered in).
$ws1.Range("H35").Value = "Martin Luther University Halle-Wittenberg"
$ws1.Range("G35").Value = "https://www.uni-halle.de/"
$ws1.Range("J35").Value = "Halle (Saale)"
$ws1.Range("D35").Value = "DiP-KS"
$ws1.Range("D35").NumberFormat = "@"
$ws1.Range("C35").Value = "https://www.dip-sachsen-anhalt.de/"
$ws1.Range("B35").Value = """"

$ws1.Range("E35").Formula = '="<a href="&B35&C35&B35&">"&D35&"</a>"'
$ws1.Range("I35").Formula = '="<a href="&B35&G35&B35&">"&H35&"</a>"'

# Hyperlinks for the newly-added URLs.
$ws1.Hyperlinks.Add($ws1.Range("G35"), "https://www.uni-halle.de/")
$ws1.Hyperlinks.Add($ws1.Range("C35"), "https://www.dip-sachsen-anhalt.de/")

# Adding a hyperlink re-stamps the cell with a fresh "Link" style; put the
# original Link style back so C35/G35 keep matching the rest of the column.
$ws1.Range("G35").Style = "Link"
$ws1.Range("C35").Style = "Link"

# Activate the People sheet/tab and select K37, scrolled so column C is
# the leftmost visible column.
$ws1.Activate()
[void]$ws1.Range("K37").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
